# Fruta / hortaliza, semanal
# Insert a new weekly price record at row 322 for "Haba" (Femacal de La Calera),
# pushing the existing rows 322:343 down to 323:344.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 322 (shifts rows 322-343 -> 323-344).
$ws.Rows(322).Insert()

# Populate the newly inserted row with the new observation.
$ws.Range("A322").Value = 3
$ws.Range("B322").Value = "Femacal de La Calera"
$ws.Range("C322").Value = "Coquimbo"
$ws.Range("D322").Value = 45265
$ws.Range("E322").Value = 5
$ws.Range("F322").Value = 100112026
$ws.Range("G322").Value = "Haba"
$ws.Range("H322").Value = "Sin especificar"
$ws.Range("I322").Value = "Primera"
$ws.Range("J322").Value = 50
$ws.Range("K322").Value = 9000
$ws.Range("L322").Value = 9000
$ws.Range("M322").Value = 9000
$ws.Range("N322").Value = "$/saco 25 kilos"
$ws.Range("O322").Value = "Provincia de Quillota"
$ws.Range("P322").Value = 360
$ws.Range("Q322").Value = 25
$ws.Range("R322").Value = "Hortaliza"
